$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 5 values to the new "custom accuracy" (2 decimal) values
$ws.Range("B5").Value  = 6.74
$ws.Range("C5").Value  = 5.32
$ws.Range("D5").Value  = 0.18
$ws.Range("E5").Value  = 14.28
$ws.Range("F5").Value  = 12.17
$ws.Range("G5").Value  = 4.95
$ws.Range("H5").Value  = 20.05
$ws.Range("I5").Value  = 7.89
$ws.Range("J5").Value  = 3.54
$ws.Range("K5").Value  = 5.18
$ws.Range("L5").Value  = 6.25
$ws.Range("M5").Value  = 6.37
$ws.Range("N5").Value  = 1.78
$ws.Range("O5").Value  = 5.11
$ws.Range("P5").Value  = 7.66
$ws.Range("Q5").Value  = 4.34
$ws.Range("R5").Value  = 0.35
$ws.Range("S5").Value  = 0.11
$ws.Range("T5").Value  = 72.55
$ws.Range("U5").Value  = 14.66
$ws.Range("V5").Value  = 4.98
$ws.Range("W5").Value  = 9.710000000000001
$ws.Range("X5").Value  = 5.04
$ws.Range("Y5").Value  = 0.68
$ws.Range("Z5").Value  = 9.75
$ws.Range("AA5").Value = 4.29
$ws.Range("AB5").Value = 3.63
$ws.Range("AC5").Value = 4.47
$ws.Range("AD5").Value = 6.12
$ws.Range("AE5").Value = 0
$ws.Range("AF5").Value = 18.35
$ws.Range("AG5").Value = 3.29
$ws.Range("AH5").Value = 5.91

# Remove row 6 entirely (data regenerated with fewer rows / "1000 data points" note)
$ws.Rows.Item(6).Delete()
